$wb = $excel.ActiveWorkbook

# The "Generate Report for handback" edit: both the zh-cn and de-de report
# sheets get a refreshed handback status, a populated "Latest Target
# File" / "Latest Handback File" pair (columns E/F) for the two source
# rows, and a new "Latest Handback DateTime" (column G) for those rows.

$sheetNames = @("zh-cn", "de-de")
$handbackDateTimes = @{ "zh-cn" = "2016-01-26 03:31:30"; "de-de" = "2016-01-26 03:31:48" }

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Collect existing hyperlink addresses for columns A and C on rows 2/3
    # so the new E/F hyperlinks can point at the very same targets.
    $aAddr2 = ""
    $cAddr2 = ""
    $aAddr3 = ""
    $cAddr3 = ""
    foreach ($h in $ws.Hyperlinks) {
        $refAddr = $h.Range.Address()
        if ($refAddr -eq '$A$2') { $aAddr2 = $h.Address() }
        if ($refAddr -eq '$C$2') { $cAddr2 = $h.Address() }
        if ($refAddr -eq '$A$3') { $aAddr3 = $h.Address() }
        if ($refAddr -eq '$C$3') { $cAddr3 = $h.Address() }
    }

    # Row 2 and row 3 text that the new hyperlinks should display (matches
    # columns A / C respectively).
    $aText2 = $ws.Range("A2").Value()
    $cText2 = $ws.Range("C2").Value()
    $aText3 = $ws.Range("A3").Value()
    $cText3 = $ws.Range("C3").Value()

    # Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # New "Latest Target File" (E) / "Latest Handback File" (F) hyperlinks,
    # mirroring the handoff file / handoff-target info already present in
    # columns A and C for each row.
    $ws.Hyperlinks.Add($ws.Range("E2"), $aAddr2, [Type]::Missing, [Type]::Missing, $aText2)
    $ws.Hyperlinks.Add($ws.Range("F2"), $cAddr2, [Type]::Missing, [Type]::Missing, $cText2)
    $ws.Hyperlinks.Add($ws.Range("E3"), $aAddr3, [Type]::Missing, [Type]::Missing, $aText3)
    $ws.Hyperlinks.Add($ws.Range("F3"), $cAddr3, [Type]::Missing, [Type]::Missing, $cText3)

    $ws.Range("E2").Style = "HyperLink"
    $ws.Range("F2").Style = "HyperLink"
    $ws.Range("E3").Style = "HyperLink"
    $ws.Range("F3").Style = "HyperLink"

    # "Latest Handback DateTime" (G) for the two source rows now reflects
    # the handback that just happened (was the epoch placeholder before).
    $dt = $handbackDateTimes[$sheetName]
    $ws.Range("G2").Value = $dt
    $ws.Range("G3").Value = $dt
}
